# Trade #19 closed at 2026-02-17 12:29:18 - unknown UNKNOWN +0.000%
#
# Appends the newly-closed trade (#19) to the "All Trades" and
# "MarketMaking" logs, and rolls the updated trade-count / win-rate
# totals into the "Summary" and "Strategy Status" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet: Total Trades 18 -> 19, Win Rate % 33.33 -> 31.58
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 19
$summary.Range("B9").Value = 31.58

# ---------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4):
# Trades 18 -> 19, Win Rate % 33.33 -> 31.58
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 19
$status.Range("G4").Value = 31.58

# ---------------------------------------------------------------
# Append the new trade row (row 20) to both "All Trades" and
# "MarketMaking" sheets - they mirror the same trade log.
# ---------------------------------------------------------------
function Add-TradeRow($ws) {
    $row = 20

    $ws.Cells.Item($row, 1).Value = 19

    # Force the date column to remain literal text ("2026-02-17")
    # instead of being auto-coerced into a date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    $ws.Cells.Item($row, 3).Value = "12:29:12"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.97
    $ws.Cells.Item($row, 7).Value = 0.97
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 99.98999999999999
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.11
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
